# Update Leve profit-calculation sheets with refreshed market data (scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 313.44
$ws.Range("I6").Value = 147
$ws.Range("K6").Value = 441
$ws.Range("M6").Value = -329

# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# Row 98
$ws.Range("H98").Value = 23589.533
$ws.Range("I98").Value = 23491.959
$ws.Range("J98").Value = 23979.834
$ws.Range("K98").Value = 23491.959
$ws.Range("L98").Value = 23979.834
$ws.Range("M98").Value = -21993.959
$ws.Range("N98").Value = -26975.834

# Row 116
$ws.Range("H116").Value = 16306586
$ws.Range("I116").Value = 28530626
$ws.Range("J116").Value = 7865
$ws.Range("K116").Value = 28530626
$ws.Range("L116").Value = 7865
$ws.Range("M116").Value = -28527184
$ws.Range("N116").Value = -14749

# Row 122
$ws.Range("H122").Value = 23589.533
$ws.Range("I122").Value = 23491.959
$ws.Range("J122").Value = 23979.834
$ws.Range("K122").Value = 70475.87699999999
$ws.Range("L122").Value = 71939.50199999999
$ws.Range("M122").Value = -68025.87699999999
$ws.Range("N122").Value = -76839.50199999999

# Row 131
$ws.Range("H131").Value = 2965.524
$ws.Range("I131").Value = 1452.8667
$ws.Range("K131").Value = 4358.6001
$ws.Range("M131").Value = 681.3999000000003

# Row 132
$ws.Range("H132").Value = 3317.1162
$ws.Range("I132").Value = 3704.8928
$ws.Range("K132").Value = 11114.6784
$ws.Range("M132").Value = -8584.678400000001

# Row 137
$ws.Range("H137").Value = 20309.666
$ws.Range("I137").Value = 25437
$ws.Range("J137").Value = 4927.6665
$ws.Range("K137").Value = 76311
$ws.Range("L137").Value = 14782.9995
$ws.Range("M137").Value = -73761
$ws.Range("N137").Value = -19882.9995

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2224.24
$ws.Range("I2").Value = 1890.45
$ws.Range("K2").Value = 1890.45
$ws.Range("M2").Value = -1777.45

# Row 6
$ws.Range("H6").Value = 4001799.5
$ws.Range("J6").Value = 2749.5
$ws.Range("L6").Value = 2749.5
$ws.Range("N6").Value = -3095.5

# Row 32
$ws.Range("H32").Value = 2544.1765
$ws.Range("I32").Value = 2257.7424
$ws.Range("K32").Value = 2257.7424
$ws.Range("M32").Value = -1970.7424

# Row 45
$ws.Range("H45").Value = 4536.4546
$ws.Range("I45").Value = 3287
$ws.Range("K45").Value = 3287
$ws.Range("M45").Value = -2910

# Row 116
$ws.Range("H116").Value = 2224.24
$ws.Range("I116").Value = 1890.45
$ws.Range("K116").Value = 1890.45
$ws.Range("M116").Value = 403.55

# Row 132
$ws.Range("H132").Value = 3710.2703
$ws.Range("I132").Value = 3697.971
$ws.Range("J132").Value = 3880
$ws.Range("K132").Value = 11093.913
$ws.Range("L132").Value = 11640
$ws.Range("M132").Value = -8563.913
$ws.Range("N132").Value = -16700

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2224.24
$ws.Range("I3").Value = 1890.45
$ws.Range("K3").Value = 1890.45
$ws.Range("M3").Value = -1776.45

# Row 22
$ws.Range("H22").Value = 292.33334
$ws.Range("I22").Value = 293.42856
$ws.Range("K22").Value = 293.42856
$ws.Range("M22").Value = -120.42856

$ws = $wb.Worksheets.Item("CRP")
# Row 9
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

# Row 31
$ws.Range("H31").Value = 4166.222
$ws.Range("I31").Value = 3456.5833
$ws.Range("J31").Value = 5585.5
$ws.Range("K31").Value = 3456.5833
$ws.Range("L31").Value = 5585.5
$ws.Range("M31").Value = -3161.5833
$ws.Range("N31").Value = -6175.5

# Row 34
$ws.Range("H34").Value = 4166.222
$ws.Range("I34").Value = 3456.5833
$ws.Range("J34").Value = 5585.5
$ws.Range("K34").Value = 3456.5833
$ws.Range("L34").Value = 5585.5
$ws.Range("M34").Value = -3254.5833
$ws.Range("N34").Value = -5989.5

# Row 134
$ws.Range("H134").Value = 6289.4165
$ws.Range("I134").Value = 3371.625
$ws.Range("K134").Value = 10114.875
$ws.Range("M134").Value = -7579.875

$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 29548.611
$ws.Range("I132").Value = 1614.5555
$ws.Range("J132").Value = 57482.668
$ws.Range("K132").Value = 14530.9995
$ws.Range("L132").Value = 517344.012
$ws.Range("M132").Value = -12000.9995
$ws.Range("N132").Value = -522404.012

$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 13667.833
$ws.Range("I22").Value = 15401.4
$ws.Range("K22").Value = 15401.4
$ws.Range("M22").Value = -14872.4

# Row 54
$ws.Range("H54").Value = 16728.572
$ws.Range("J54").Value = 16728.572
$ws.Range("L54").Value = 16728.572
$ws.Range("N54").Value = -17508.572

# Row 57
$ws.Range("H57").Value = 33331.5
$ws.Range("J57").Value = 39994.5
$ws.Range("L57").Value = 39994.5
$ws.Range("N57").Value = -41634.5

# Row 59
$ws.Range("H59").Value = 6820
$ws.Range("I59").Value = 7333.3335
$ws.Range("J59").Value = 6050
$ws.Range("K59").Value = 7333.3335
$ws.Range("L59").Value = 6050
$ws.Range("M59").Value = -6750.3335
$ws.Range("N59").Value = -7216

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 2778.6191
$ws.Range("I82").Value = 2579.6667
$ws.Range("K82").Value = 2579.6667
$ws.Range("M82").Value = -2218.6667

# Row 85
$ws.Range("H85").Value = 2778.6191
$ws.Range("I85").Value = 2579.6667
$ws.Range("K85").Value = 2579.6667
$ws.Range("M85").Value = -1331.6667

# Row 122
$ws.Range("H122").Value = 6227.5137
$ws.Range("I122").Value = 4558.1816
$ws.Range("J122").Value = 19999.5
$ws.Range("K122").Value = 13674.5448
$ws.Range("L122").Value = 59998.5
$ws.Range("M122").Value = -11224.5448
$ws.Range("N122").Value = -64898.5

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 975.36365
$ws.Range("I4").Value = 85.8
$ws.Range("J4").Value = 1716.6666
$ws.Range("K4").Value = 85.8
$ws.Range("L4").Value = 1716.6666
$ws.Range("M4").Value = 27.2
$ws.Range("N4").Value = -1942.6666

# Row 81
$ws.Range("H81").Value = 10358.667
$ws.Range("J81").Value = 3757.6
$ws.Range("L81").Value = 7515.2
$ws.Range("N81").Value = -9637.200000000001

# Row 84
$ws.Range("H84").Value = 10358.667
$ws.Range("J84").Value = 3757.6
$ws.Range("L84").Value = 37576
$ws.Range("N84").Value = -48184

# Row 96
$ws.Range("H96").Value = 2246.6
$ws.Range("I96").Value = 2228.5
$ws.Range("J96").Value = 2258.6667
$ws.Range("K96").Value = 2228.5
$ws.Range("L96").Value = 2258.6667
$ws.Range("M96").Value = -855.5
$ws.Range("N96").Value = -5004.6667

# Row 136
$ws.Range("H136").Value = 1862.6578
$ws.Range("I136").Value = 1218.9354
$ws.Range("K136").Value = 3656.8062
$ws.Range("M136").Value = -1106.8062
